# Weekly refresh of Fruta / hortaliza data:
# the rows keep the same static descriptive columns (A,B,C,E,F,G,H,I,N,O,Q,R)
# but the date (D) together with the volume/price columns (J,K,L,M,P) are
# re-sorted to reflect the latest weekly pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => Fecha (D, Excel date serial number), Volumen (J), Precio minimo (K),
#        Precio maximo (L), Precio promedio ponderado (M), Precio $/Kg (P)
$rows = @(
    @{ Row = 2;  D = 44473; J = 140; K = 1600; L = 1600; M = 1600; P = 1600 },
    @{ Row = 3;  D = 44497; J = 50;  K = 2200; L = 2200; M = 2200; P = 2200 },
    @{ Row = 4;  D = 44483; J = 50;  K = 2200; L = 2200; M = 2200; P = 2200 },
    @{ Row = 5;  D = 44476; J = 30;  K = 2200; L = 2200; M = 2200; P = 2200 },
    @{ Row = 6;  D = 44487; J = 50;  K = 2200; L = 2200; M = 2200; P = 2200 },
    @{ Row = 7;  D = 44484; J = 40;  K = 2200; L = 2200; M = 2200; P = 2200 },
    @{ Row = 8;  D = 44452; J = 120; K = 2300; L = 2300; M = 2300; P = 2300 },
    @{ Row = 9;  D = 44203; J = 30;  K = 2000; L = 2000; M = 2000; P = 2000 },
    @{ Row = 10; D = 44447; J = 75;  K = 2200; L = 2200; M = 2200; P = 2200 },
    @{ Row = 11; D = 44496; J = 40;  K = 2200; L = 2200; M = 2200; P = 2200 },
    @{ Row = 12; D = 44453; J = 20;  K = 2300; L = 2300; M = 2300; P = 2300 },
    @{ Row = 13; D = 44474; J = 20;  K = 1600; L = 1600; M = 1600; P = 1600 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value = $r.D        # D - Fecha
    $ws.Cells.Item($r.Row, 10).Value = $r.J       # J - Volumen
    $ws.Cells.Item($r.Row, 11).Value = $r.K       # K - Precio minimo
    $ws.Cells.Item($r.Row, 12).Value = $r.L       # L - Precio maximo
    $ws.Cells.Item($r.Row, 13).Value = $r.M       # M - Precio promedio ponderado
    $ws.Cells.Item($r.Row, 16).Value = $r.P       # P - Precio $/Kg
}
